$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$aw = $excel.ActiveWindow
Write-Host "WindowHeight prop?" $aw.WindowHeight
